$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the header style used by G1 ("sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new values in H2 and H3 (plain numeric cells, same as rest of data rows)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
